$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.337343811988831
$ws.Range("B1").Value = 1.547056913375854
$ws.Range("C1").Value = 4.000873565673828
$ws.Range("D1").Value = 3.208665370941162
$ws.Range("E1").Value = 1.095940947532654
